$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Metadata")

# Bump "Version" value: 1.0.0 -> 1.1.0
$ws1.Range("B3").Value = "1.1.0"

# Update "Date" value
$ws1.Range("B8").Value = "2023-07-10T23:08:03+02:00"

# Re-apply (add) the alignment/wrap formatting to the already-wrapped ranges so
# the cell style definitions carry applyAlignment="true" (matches the edited
# workbook's cellXfs, which gained that attribute on both custom styles).
$ws1.Range("A1:B1").WrapText = $true
$ws1.Range("A2:B14").WrapText = $true

$ws2 = $wb.Worksheets.Item("Include from CareSocialCodes")
$ws2.Range("A1:C1").WrapText = $true
$ws2.Range("A2:B4").WrapText = $true
$ws2.Range("C2").WrapText = $true
